$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06747562906149941
$ws.Range("C2").Value = 0.9987673800532109
$ws.Range("D2").Value = 0.2063647345469068
$ws.Range("G2").Value = 0.2675185060000028
$ws.Range("H2").Value = 0.987

$ws.Range("B3").Value = 0.1216235692055888
$ws.Range("C3").Value = 0.9910140076357281
$ws.Range("D3").Value = 0.2726675398736592
$ws.Range("G3").Value = 0.2675185060000028
$ws.Range("H3").Value = 0.987
